$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 181.4944075
$ws.Range("H2").Value = 362.988815
$ws.Range("I2").Value = 0.2239486468210351
$ws.Range("J2").Value = 0.1654349085470023
$ws.Range("Q2").Value = 57.06740754649667
$ws.Range("R2").Value = 342.4044452789801
$ws.Range("S2").Value = 0.2239486468210351
$ws.Range("T2").Value = 0.1654349085470023

# Row 3
$ws.Range("I3").Value = 0.07700606288633029
$ws.Range("J3").Value = 0.08532865336765341
$ws.Range("S3").Value = 0.07700606288633029
$ws.Range("T3").Value = 0.08532865336765341

# Row 4
$ws.Range("G4").Value = 171.9980316666667
$ws.Range("H4").Value = 515.994095
$ws.Range("I4").Value = 0.2122309275432167
$ws.Range("J4").Value = 0.235168226649403
$ws.Range("Q4").Value = 54.08145576230445
$ws.Range("R4").Value = 486.73310186074
$ws.Range("S4").Value = 0.2122309275432167
$ws.Range("T4").Value = 0.235168226649403

# Row 5
$ws.Range("G5").Value = 55.64279550000001
$ws.Range("H5").Value = 111.285591
$ws.Range("I5").Value = 0.06865847234198982
$ws.Range("J5").Value = 0.05071925307032974
$ws.Range("Q5").Value = 17.495801284262
$ws.Range("R5").Value = 104.974807705572
$ws.Range("S5").Value = 0.06865847234198982
$ws.Range("T5").Value = 0.05071925307032974

# Row 6
$ws.Range("G6").Value = 203.386317
$ws.Range("H6").Value = 610.158951
$ws.Range("I6").Value = 0.250961399315095
$ws.Range("J6").Value = 0.2780845747487284
$ws.Range("Q6").Value = 63.950895245188
$ws.Range("R6").Value = 575.5580572066921
$ws.Range("S6").Value = 0.250961399315095
$ws.Range("T6").Value = 0.2780845747487284

# Row 7
$ws.Range("G7").Value = 135.4992116666667
$ws.Range("H7").Value = 406.497635
$ws.Range("I7").Value = 0.167194491092333
$ws.Range("J7").Value = 0.1852643836168829
$ws.Range("Q7").Value = 42.60510745715778
$ws.Range("R7").Value = 383.44596711442
$ws.Range("S7").Value = 0.167194491092333
$ws.Range("T7").Value = 0.1852643836168829
